$d = $word.ActiveDocument

# "Version 1." -> "Version 2." but re-shaped at the run level:
#   - "Version" splits into two runs: "Versi" + "on"
#   - " 1." becomes " 2" (drop the trailing period from that run)
#   - a brand-new trailing run holding just "." is appended after the
#     _GoBack bookmark

# 1) Split "Version" into "Versi" / "on" by dropping a temporary,
#    zero-length bookmark at the boundary (char offset 5) and immediately
#    removing it. Word always breaks the run at a bookmark's position, and
#    deleting the bookmark afterwards leaves the break in place without
#    adding any direct formatting residue.
$splitPoint = $d.Range(5, 5)
$d.Bookmarks.Add("tmpSplit", $splitPoint)
$d.Bookmarks("tmpSplit").Delete()

# 2) " 1." -> " 2" (the run keeps its leading space, loses the period)
$d.Content.Find.Execute(" 1.", $true, $false, $false, $false, $false, $true, 1, $false, " 2", 2)

# 3) Re-append the period as its own trailing run, after the bookmark.
$endRange = $d.Range(10, 10)
$endRange.InsertAfter(".")
